$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1 (Week 3 paragraph): insert a comma after "...for our project"
# "...help for our project as we got..." -> "...help for our project, as we got..."
# ---------------------------------------------------------------------------
$week3 = $d.Paragraphs.Item(6).Range.Duplicate
$week3.Find.Execute("help for our project", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$week3.Collapse(0)
$week3.InsertAfter(",")

# ---------------------------------------------------------------------------
# Change 2 (Week 4 paragraph): append the "This week, Shweta created..." text
# after the existing tab character in paragraph 8.
# ---------------------------------------------------------------------------
$week4Body = $d.Paragraphs.Item(8).Range
$week4Body.Collapse(0)
$week4Body.InsertAfter("This week")
$week4Body = $d.Paragraphs.Item(8).Range
$week4Body.Collapse(0)
$week4Body.InsertAfter(",")
$week4Body = $d.Paragraphs.Item(8).Range
$week4Body.Collapse(0)
$week4Body.InsertAfter(" Shweta created a react app in our repo and added some code for responsive headers with the help of ‘react-bootstrap’ module. I modified our project folder structure and added some routing functionality between home-page, hike-trails page, weather-page. ")

# ---------------------------------------------------------------------------
# Change 3: append Week 5 .. Week 10 header/TODO paragraphs.
# ---------------------------------------------------------------------------
function Add-WeekHeaderParagraph($afterParaIndex, $weekNumber) {
    $tail = $d.Paragraphs.Item($afterParaIndex).Range
    $tail.Collapse(0)
    $tail.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($afterParaIndex + 1).Range
    $newPara.Text = "Week "
    $newPara = $d.Paragraphs.Item($afterParaIndex + 1).Range
    $newPara.Collapse(0)
    $newPara.InsertAfter("$weekNumber")
    $newPara = $d.Paragraphs.Item($afterParaIndex + 1).Range
    $newPara.Collapse(0)
    $newPara.InsertAfter(" –")
    $whole = $d.Paragraphs.Item($afterParaIndex + 1).Range
    $whole.Bold = 1
    $whole.BoldBi = 1
}

function Add-TodoParagraph($afterParaIndex, $leadingSpace) {
    $tail = $d.Paragraphs.Item($afterParaIndex).Range
    $tail.Collapse(0)
    $tail.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($afterParaIndex + 1).Range
    if ($leadingSpace) {
        $newPara.Text = " "
    }
    $run1 = $d.Paragraphs.Item($afterParaIndex + 1).Range
    $run1.Collapse(1)
    $run1.InsertBefore([char]9)
    $newPara = $d.Paragraphs.Item($afterParaIndex + 1).Range
    $newPara.Collapse(0)
    $newPara.InsertAfter("<!--  TODO  -->")
}

$idx = 8
Add-WeekHeaderParagraph $idx 5
$idx = $idx + 1
Add-TodoParagraph $idx $true
$idx = $idx + 1
Add-WeekHeaderParagraph $idx 6
$idx = $idx + 1
Add-TodoParagraph $idx $false
$idx = $idx + 1
Add-WeekHeaderParagraph $idx 7
$idx = $idx + 1
Add-TodoParagraph $idx $false
$idx = $idx + 1
Add-WeekHeaderParagraph $idx 8
$idx = $idx + 1
Add-TodoParagraph $idx $false
$idx = $idx + 1
Add-WeekHeaderParagraph $idx 9
$idx = $idx + 1
Add-TodoParagraph $idx $false
$idx = $idx + 1
Add-WeekHeaderParagraph $idx 10
$idx = $idx + 1
Add-TodoParagraph $idx $false

Write-Output "Paragraph count: $($d.Paragraphs.Count)"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    Write-Output "Para $i : [$($p.Range.Text)]"
}
